# Updates the cryptos list (prices + 1h volume deltas) per the
# Sun Oct 20 07:37:50 UTC 2024 GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text fields (coin name / link / volume%) are never numeric-
# looking, so a direct .Value assignment keeps them as text.
$ws.Range("E2").Value = "  +0.06%  "
$ws.Range("E3").Value = "  +0.24%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("E6").Value = "  +2.91%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  -0.91%  "
$ws.Range("E9").Value = "  -1.02%  "
$ws.Range("E10").Value = "  -1.04%  "
$ws.Range("E11").Value = "  +0.51%  "
$ws.Range("E12").Value = "  -0.75%  "
$ws.Range("E13").Value = "  -0.61%  "
$ws.Range("E14").Value = "  +0.33%  "
$ws.Range("E15").Value = "  -2.92%  "
$ws.Range("E16").Value = "  +0.12%  "
$ws.Range("E17").Value = "  +0.41%  "
$ws.Range("E18").Value = "  -0.27%  "
$ws.Range("E19").Value = "  -1.03%  "
$ws.Range("E20").Value = "  -1.31%  "
$ws.Range("E22").Value = "  -2.62%  "
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("E24").Value = "  -0.12%  "
$ws.Range("E25").Value = "  -0.08%  "
$ws.Range("E26").Value = "  -0.74%  "
$ws.Range("E28").Value = "  -2.63%  "
$ws.Range("E29").Value = "  -0.18%  "
$ws.Range("E30").Value = "  -2.22%  "
$ws.Range("E31").Value = "  -1.68%  "
$ws.Range("E32").Value = "  -1.96%  "
$ws.Range("E33").Value = "  +1.04%  "
$ws.Range("E34").Value = "  +3.33%  "
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("E36").Value = "  -1.85%  "
$ws.Range("B37").Value = "Monero"
$ws.Range("C37").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("E37").Value = "  -0.81%  "
$ws.Range("B38").Value = "EthereumClassic"
$ws.Range("C38").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("E38").Value = "  +1.56%  "
$ws.Range("E40").Value = "  -1.51%  "
$ws.Range("E41").Value = "  -1.55%  "
$ws.Range("E42").Value = "  -1.45%  "
$ws.Range("E43").Value = "  -5.17%  "
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("E45").Value = "  +0.67%  "
$ws.Range("E46").Value = "  +0.81%  "
$ws.Range("E47").Value = "  +0.40%  "
$ws.Range("E48").Value = "  -1.49%  "
$ws.Range("E49").Value = "  -2.02%  "
$ws.Range("E50").Value = "  +1.21%  "
$ws.Range("E51").Value = "  -0.23%  "

# Price fields (column D) are stored as text in the source sheet, but
# several look like plain numbers (e.g. "5.29") and Excel auto-converts
# those to the Number type on a bare .Value assignment. Force the Text
# number format while writing, then clear formatting again afterwards
# so the cell keeps its original (unstyled) appearance while the stored
# value stays a genuine string, matching the workbook's original layout.
$priceCells = @{
    "D2" = "68.415.21"
    "D3" = "2.648.11"
    "D5" = "597.94"
    "D6" = "159.05"
    "D11" = "5.29"
    "D13" = "28.01"
    "D14" = "3.132.00"
    "D15" = "0.0000188"
    "D16" = "68.317.58"
    "D17" = "2.653.10"
    "D19" = "360.38"
    "D20" = "7.41"
    "D22" = "4.77"
    "D24" = "74.57"
    "D25" = "0.999"
    "D26" = "9.77"
    "D27" = "2.785.57"
    "D29" = "0.999"
    "D30" = "561.36"
    "D32" = "1.40"
    "D34" = "1.65"
    "D37" = "159.62"
    "D38" = "19.70"
    "D45" = "157.44"
    "D47" = "22.02"
    "D49" = "0.0772"
    "D51" = "0.615"
}

foreach ($addr in $priceCells.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $priceCells[$addr]
    $cell.ClearFormats()
}

Write-Output "Applied cryptos update"
